# 8.3.1.2.xlsx — add the 2020 column (K) to the indicator table and
# move the active-cell selection, matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K3: "2020", styled like the other year headers (bold,
# Times New Roman 10pt, top+bottom medium border) ---
$ws.Range("I3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2020
$ws.Range("K3").Font.Name = "Times New Roman"
$ws.Range("K3").Font.Size = 10
$ws.Range("K3").Font.Bold = $true

# --- Data cells K4/K5: new 2020 figures, right aligned, thousands-style
# decimal format, smaller Kyrghyz Times font, matching row borders ---
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2.8218550629805335
$ws.Range("K4").NumberFormat = "#,##0.0"
$ws.Range("K4").Font.Name = "Kyrghyz Times"
$ws.Range("K4").Font.Size = 9
$ws.Range("K4").HorizontalAlignment = -4152

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 1.3005071159823327
$ws.Range("K5").NumberFormat = "#,##0.0"
$ws.Range("K5").Font.Name = "Kyrghyz Times"
$ws.Range("K5").Font.Size = 9
$ws.Range("K5").HorizontalAlignment = -4152

# --- Move the active selection as recorded in the saved view state ---
$null = $ws.Range("L8").Select()
